$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "ceiling type" cell J3: was the shared string "Roof_BC1", now a brand-new
# shared string "Floor_BC1" (adds a new entry to sharedStrings.xml).
$ws.Range("J3").Value = "Floor_BC1"

# "PV tilt" cell AI3: 0 -> 30
$ws.Range("AI3").Value = 30

# Move the on-screen selection/active cell to H7 (view scrolled back towards
# the start of the sheet as part of this edit).
$ws.Activate()
$ws.Range("H7").Select()
